$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17 ---
$ws.Range("A17").Value = 43742
$ws.Range("B17").Value = 0.91666666666666663
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = "Node js express 강의 수강"
$ws.Range("F17").Characters(19, 3).Font.Name = "돋움"
$ws.Range("F17").Characters(19, 3).Font.Size = 10

# --- Row 18 ---
$ws.Range("A18").Value = 43744
$ws.Range("B18").Value = 0.79166666666666663
$ws.Range("C18").Value = 0.041666666666666664
$ws.Range("D18").Value = 90
$ws.Range("E18").Value = 270
$ws.Range("F18").Value = "프로토타입 express framework 적용 / refactoring"

# --- Row 19 ---
$ws.Range("A19").Value = 43746
$ws.Range("B19").Value = 0.5
$ws.Range("C19").Value = 0.54166666666666663
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = "프로토타입 express framework 적용 / refactoring"
$ws.Range("F18").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 13

# --- Selection / view state ---
$ws.Range("F1048551").Select()
